# Dochadzka.xlsx - Autentifikacia Admin / Uctovnicka/User
# Applies:
#  1. Row 134 (Dochadzka sheet): Cas 11:30 -> 11:40
#  2. Row 264: admin entry 08:02/Odchod -> 08:10/Prichod
#  3. Old row 268 (admin,16,11:00,Prichod) removed (rows below shift up)
#  4. 12 new attendance rows appended at the end (rows 279-290)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dochadzka")

# 1. Fix row 134 (jan cyprich, apríl, 10, Odchod) - Cas 11:30 -> 11:40
$ws.Cells.Item(134, 4).Value = "11:40"

# 2. Row 264: change the time + type (admin, apríl, 16)
$ws.Cells.Item(264, 4).Value = "08:10"
$ws.Cells.Item(264, 5).Value = "Prichod"

# 3. Remove the old row 268 (admin, apríl, 16, 11:00, Prichod) - everything below moves up one row
$ws.Rows.Item(268).Delete()

# 4. Append the new attendance rows at the bottom of the table (now rows 279-290)
$newRows = @(
    @{ A = "admin";           C = "apríl"; row = 16; D = "14:53"; E = "Odchod" },
    @{ A = "admin";           C = "apríl"; row = 16; D = "06:00"; E = "Odchod" },
    @{ A = "admin";           C = "apríl"; row = 16; D = "06:14"; E = "Odchod" },
    @{ A = "jozef florek";    C = "apríl"; row = 17; D = "12:08"; E = "Odchod" },
    @{ A = "jan cyprich";     C = "apríl"; row = 10; D = "16:00"; E = "Odchod" },
    @{ A = "milada";          C = "apríl"; row = 17; D = "14:55"; E = "Prichod" },
    @{ A = "Viliam Jasurek";  C = "apríl"; row = 17; D = "15:59"; E = "Odchod" },
    @{ A = "jan cyprich";     C = "apríl"; row = 17; D = "15:59"; E = "Prichod" },
    @{ A = "jozef florek";    C = "apríl"; row = 17; D = "16:00"; E = "Prichod" },
    @{ A = "Martin Straka";   C = "apríl"; row = 17; D = "16:01"; E = "Prichod" },
    @{ A = "jozef florek";    C = "apríl"; row = 17; D = "16:02"; E = "Odchod" },
    @{ A = "Augustin Straka"; C = "apríl"; row = 17; D = "16:02"; E = "Prichod" }
)

$startRow = 279
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $item = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $item.A
    $ws.Cells.Item($r, 2).Value = $item.C
    $ws.Cells.Item($r, 3).Value = $item.row
    $ws.Cells.Item($r, 4).Value = $item.D
    $ws.Cells.Item($r, 5).Value = $item.E
}
